# Insert a new paragraph ("And yes I know about circle referencing issues...")
# right after the Memory Model section (after "All 'heap' memory access...")
# and before the "Threads" heading, carrying the hidden _GoBack bookmark with
# it (moved off the end of the "marshaling" paragraph, which is where it
# used to sit).

$d = $word.ActiveDocument

# Paragraph 7 is "All "heap" memory access is asynchronous. ..." - the last
# paragraph of the Memory Model section, right before the "Threads" heading.
$memoryModelLastPara = $d.Paragraphs.Item(7)

# Split a new, empty paragraph in right after it.
$newParaRange = $memoryModelLastPara.Range.InsertParagraphAfter()

# Re-fetch the document/paragraph collection since the content shifted.
$d2 = $word.ActiveDocument
$newPara = $d2.Paragraphs.Item(8)

# Curly single quotes (U+2019), matching the rest of the document's style.
$rsquo = [char]0x2019

# Set the new paragraph's text. A single placeholder character "X" is
# temporarily appended at the very end - this gives us a safe (non
# paragraph-boundary) position at which to (re)plant the collapsed
# "_GoBack" bookmark; collapsed bookmarks placed exactly on the last
# paragraph-mark boundary land in the wrong spot, but one character before
# that boundary is reliable. We delete the placeholder afterwards and the
# bookmark stays put, collapsed right after the real text.
$newPara.Range.Text = "And yes I know about circle referencing issues. Currently I honestly don" + $rsquo + "t have a solution and I am just going with hoping it won" + $rsquo + "t cause problems.X"

$d3 = $word.ActiveDocument
$newPara2 = $d3.Paragraphs.Item(8)

# Position right before the trailing "X" placeholder.
$bookmarkPos = $newPara2.Range.End - 2
$bookmarkRange = $d3.Range($bookmarkPos, $bookmarkPos)

# "_GoBack" is Word's special, singleton "last edit" bookmark - re-adding it
# under this name moves the existing one rather than creating a duplicate,
# so it disappears from the end of the "marshaling" paragraph automatically.
$d3.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the temporary placeholder character now that the bookmark is
# anchored in the right place.
$d4 = $word.ActiveDocument
$placeholderRange = $d4.Range($bookmarkPos, $bookmarkPos + 1)
$placeholderRange.Delete()
